# [Soubarnika] Added negative test cases
# Adds two new rows of test data to the Login_TestData sheet:
#   Row 3: a phone-number style value in column A only (negative case, no password)
#   Row 4: the same phone-number style value in column A plus a new password
#          value "muthu96" in column B
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 9629459258

$ws.Range("A4").Value = 9629459258
$ws.Range("B4").Value = "muthu96"

# Match the author's final selection/active cell (B4)
$ws.Range("B4").Select() | Out-Null
